$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:EUROTEXIND"
$ws.Range("C2").Value = "NSE:ACI"
$ws.Range("E2").Value = "NSE:ADANIENT"
$ws.Range("F2").Value = ""
$ws.Range("B3").Value = "NSE:HOMEFIRST"
$ws.Range("C3").Value = "NSE:AMNPLST"
$ws.Range("E3").Value = "NSE:FORTIS"
$ws.Range("F3").Value = ""
$ws.Range("B4").Value = "NSE:PASUPTAC"
$ws.Range("C4").Value = "NSE:APOLLOPIPE"
$ws.Range("E4").Value = "NSE:JIOFIN"
$ws.Range("F4").Value = ""
$ws.Range("B5").Value = "NSE:RHL"
$ws.Range("C5").Value = "NSE:AXISHCETF"
$ws.Range("E5").Value = "NSE:NCC"
$ws.Range("F5").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "NSE:AXISNIFTY"
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "NSE:BIGBLOC"
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "NSE:CHAMBLFERT"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "NSE:CLEDUCATE"
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "NSE:DBREALTY"
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "NSE:DCMNVL"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "NSE:DJML"
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "NSE:DPWIRES"
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "NSE:DREDGECORP"
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "NSE:FDC"
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "NSE:GEOJITFSL"
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "NSE:GREENPLY"
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "NSE:GSS"
$ws.Range("E18").Value = ""
$ws.Range("F18").Value = ""
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = "NSE:GTLINFRA"
$ws.Range("E19").Value = ""
$ws.Range("F19").Value = ""
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = "NSE:HEALTHY"
$ws.Range("E20").Value = ""
$ws.Range("F20").Value = ""
$ws.Range("B21").Value = ""
$ws.Range("C21").Value = "NSE:HIKAL"
$ws.Range("E21").Value = ""
$ws.Range("F21").Value = ""
$ws.Range("B22").Value = ""
$ws.Range("C22").Value = "NSE:IGL"
$ws.Range("E22").Value = ""
$ws.Range("F22").Value = ""
$ws.Range("B23").Value = ""
$ws.Range("C23").Value = "NSE:KDDL"
$ws.Range("E23").Value = ""
$ws.Range("F23").Value = ""
$ws.Range("B24").Value = ""
$ws.Range("C24").Value = "NSE:KHANDSE"
$ws.Range("E24").Value = ""
$ws.Range("F24").Value = ""
$ws.Range("B25").Value = ""
$ws.Range("C25").Value = "NSE:LALPATHLAB"
$ws.Range("E25").Value = ""
$ws.Range("F25").Value = ""
$ws.Range("B26").Value = ""
$ws.Range("C26").Value = "NSE:LOVABLE"
$ws.Range("E26").Value = ""
$ws.Range("F26").Value = ""
$ws.Range("B27").Value = ""
$ws.Range("C27").Value = "NSE:MAHSEAMLES"
$ws.Range("E27").Value = ""
$ws.Range("F27").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("C28").Value = "NSE:MALUPAPER"
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = "NSE:MASTEK"
$ws.Range("E29").Value = ""
$ws.Range("F29").Value = ""
$ws.Range("B30").Value = ""
$ws.Range("C30").Value = "NSE:OSWALAGRO"
$ws.Range("E30").Value = ""
$ws.Range("F30").Value = ""
$ws.Range("B31").Value = ""
$ws.Range("C31").Value = "NSE:PODDARMENT"
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = ""
$ws.Range("B32").Value = ""
$ws.Range("C32").Value = "NSE:RADIANTCMS"
$ws.Range("E32").Value = ""
$ws.Range("F32").Value = ""
$ws.Range("B33").Value = ""
$ws.Range("C33").Value = "NSE:RAMCOIND"
$ws.Range("E33").Value = ""
$ws.Range("F33").Value = ""
$ws.Range("B34").Value = ""
$ws.Range("C34").Value = "NSE:ROLEXRINGS"
$ws.Range("E34").Value = ""
$ws.Range("F34").Value = ""
$ws.Range("B35").Value = ""
$ws.Range("C35").Value = "NSE:ROML"
$ws.Range("E35").Value = ""
$ws.Range("F35").Value = ""
$ws.Range("B36").Value = ""
$ws.Range("C36").Value = "NSE:RTNINDIA"
$ws.Range("E36").Value = ""
$ws.Range("F36").Value = ""

# Remove now-obsolete rows 37:49 (table shrinks from 49 to 36 rows)
$ws.Range("A37:F49").EntireRow.Delete()
